$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BS_DATA_VUE")

# --- Row 3 fix-ups: fill in the previously-empty B3, and correct D3 ---
$ws.Range("B3").Value = "As a user, I want to manage my toDo list"
$ws.Range("D3").Value = "NIL"

# --- Move the old row 4 (VUE_003 -> renamed to VUE_004) down to row 5 ---
# Copy format+value first, then overwrite the text that changed.
$ws.Range("A4:E4").Copy()
$ws.Range("A5:E5").PasteSpecial(-4122)
$ws.Range("A5").Value = "VUE_004_UserShouldNotBeAbleToWriteMoreThan20CharsLong"
$ws.Range("B5").Value = "As a user, I should expect to see error message if I create more than 25 characters long"
$ws.Range("C5").Value = "Negative"
$ws.Range("D5").Value = "NIL"
$ws.Range("E5").Value = "Wash the cat | asdfadsnoarenfahoiegrfjsaigtselfdsnbsukjfedasoidvjg"

# --- Populate the new row 4 with the Chinese test case, using row 2's formatting ---
$ws.Range("A2:E2").Copy()
$ws.Range("A4:E4").PasteSpecial(-4122)
$ws.Range("A4").Value = "VUE_003_UserShouldBeAbleToPopulateToDosWithChinese"
$ws.Range("B4").Value = "As a user, I want to manage my toDo list"
$ws.Range("C4").Value = "Positive"
$ws.Range("D4").Value = "NIL"
$ws.Range("E4").Value = "我需要在早上6.30醒來 | 我需要刷牙"
$ws.Range("E4").WrapText = $true

$null = $ws.Range("E4").Select()
